$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain text (not auto-converted to numbers),
# matching the original inline-string cell type.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D18").NumberFormat = "@"
$ws.Range("D20:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '24.669.13'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '1.700.43'
$ws.Range("E3").Value = '  +0.95%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '309.12'
$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("D6").Value = '0.9976'
$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").Value = '0.3736'
$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("D8").Value = '48.92'
$ws.Range("E8").Value = '  +2.24%  '

$ws.Range("D9").Value = '0.3437'
$ws.Range("E9").Value = '  -0.77%  '

$ws.Range("D10").Value = '1.182'
$ws.Range("E10").Value = '  -1.01%  '

$ws.Range("D11").Value = '0.07458'
$ws.Range("E11").Value = '  +1.75%  '

$ws.Range("D12").Value = '0.9971'
$ws.Range("E12").Value = '  -0.44%  '

$ws.Range("D13").Value = '20.89'
$ws.Range("E13").Value = '  +1.62%  '

$ws.Range("D14").Value = '6.229'
$ws.Range("E14").Value = '  +1.54%  '

$ws.Range("D15").Value = '6.925'
$ws.Range("E15").Value = '  +1.79%  '

$ws.Range("D16").Value = '1.705.73'
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").Value = '0.00001120'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").Value = '0.06704'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = '83.65'
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("D21").Value = '17.07'
$ws.Range("E21").Value = '  +3.27%  '

$ws.Range("D22").Value = '6.326'
$ws.Range("E22").Value = '  +3.10%  '

$ws.Range("D23").Value = '13.21'
$ws.Range("E23").Value = '  +9.37%  '

$ws.Range("D24").Value = '24.672.75'
$ws.Range("E24").Value = '  +1.96%  '

$ws.Range("D25").Value = '2.408'
$ws.Range("E25").Value = '  -0.33%  '

$ws.Range("D26").Value = '2.756'
$ws.Range("E26").Value = '  +2.30%  '

$ws.Range("D27").Value = '20.09'
$ws.Range("E27").Value = '  +1.90%  '

$ws.Range("D28").Value = '150.17'
$ws.Range("E28").Value = '  -1.62%  '

$ws.Range("D29").Value = '131.04'
$ws.Range("E29").Value = '  +3.13%  '

$ws.Range("D30").Value = '1.893.70'
$ws.Range("E30").Value = '  +1.28%  '

$ws.Range("D31").Value = '1.186'
$ws.Range("E31").Value = '  +19.34%  '

$ws.Range("D32").Value = '6.762'
$ws.Range("E32").Value = '  +4.18%  '

$ws.Range("D33").Value = '4.189'
$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.779'
$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08789'
$ws.Range("E35").Value = '  +3.53%  '

$ws.Range("D36").Value = '13.64'
$ws.Range("E36").Value = '  +8.77%  '

$ws.Range("D37").Value = '5.516'
$ws.Range("E37").Value = '  +2.05%  '

$ws.Range("D38").Value = '0.06503'
$ws.Range("E38").Value = '  +0.04%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02375'
$ws.Range("E39").Value = '  +0.44%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.906'
$ws.Range("E40").Value = '  -0.51%  '

$ws.Range("D41").Value = '0.2219'
$ws.Range("E41").Value = '  +3.54%  '

$ws.Range("D42").Value = '1.273'
$ws.Range("E42").Value = '  -0.61%  '

$ws.Range("D43").Value = '0.6399'
$ws.Range("E43").Value = '  +2.87%  '

$ws.Range("D44").Value = '0.9970'
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").Value = '13.84'
$ws.Range("E45").Value = '  +3.61%  '

$ws.Range("D46").Value = '0.6083'
$ws.Range("E46").Value = '  +1.59%  '

$ws.Range("D47").Value = '3.810'
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").Value = '2.112'
$ws.Range("E48").Value = '  +3.42%  '

$ws.Range("D49").Value = '129.20'
$ws.Range("E49").Value = '  +0.96%  '

$ws.Range("E50").Value = '  +0.98%  '

$ws.Range("D51").Value = '78.94'
$ws.Range("E51").Value = '  +3.19%  '
